$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 100545.2
$ws.Range("I41").Value = 181
$ws.Range("K41").Value = 181
$ws.Range("M41").Value = 259
$ws.Range("H80").Value = 608.8788
$ws.Range("I80").Value = 312.30768
$ws.Range("J80").Value = 801.65
$ws.Range("K80").Value = 936.92304
$ws.Range("L80").Value = 2404.95
$ws.Range("M80").Value = 61.07695999999999
$ws.Range("N80").Value = -4400.95
$ws.Range("H83").Value = 608.8788
$ws.Range("I83").Value = 312.30768
$ws.Range("J83").Value = 801.65
$ws.Range("K83").Value = 2810.76912
$ws.Range("L83").Value = 7214.849999999999
$ws.Range("M83").Value = 2181.23088
$ws.Range("N83").Value = -17198.85
$ws.Range("H98").Value = 1298.08
$ws.Range("I98").Value = 1107.2325
$ws.Range("K98").Value = 1107.2325
$ws.Range("M98").Value = 390.7674999999999
$ws.Range("H100").Value = 6867.8667
$ws.Range("I100").Value = 3399.2856
$ws.Range("K100").Value = 3399.2856
$ws.Range("M100").Value = -2858.2856
$ws.Range("H122").Value = 1298.08
$ws.Range("I122").Value = 1107.2325
$ws.Range("K122").Value = 3321.6975
$ws.Range("M122").Value = -871.6975000000002
$ws.Range("H135").Value = 1598.7037
$ws.Range("I135").Value = 1621.7307
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 14595.5763
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -12060.5763
$ws.Range("N135").Value = -14070
$ws.Range("H138").Value = 3177.2856
$ws.Range("J138").Value = 4516.15
$ws.Range("L138").Value = 13548.45
$ws.Range("N138").Value = -23828.45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5956109.5
$ws.Range("I74").Value = 7815122.5
$ws.Range("J74").Value = 7268
$ws.Range("K74").Value = 7815122.5
$ws.Range("L74").Value = 7268
$ws.Range("M74").Value = -7814248.5
$ws.Range("N74").Value = -9016
$ws.Range("H77").Value = 5956109.5
$ws.Range("I77").Value = 7815122.5
$ws.Range("J77").Value = 7268
$ws.Range("K77").Value = 39075612.5
$ws.Range("L77").Value = 36340
$ws.Range("M77").Value = -39071244.5
$ws.Range("N77").Value = -45076
$ws.Range("H102").Value = 3501.0476
$ws.Range("I102").Value = 2640.0557
$ws.Range("K102").Value = 2640.0557
$ws.Range("M102").Value = -1018.0557
$ws.Range("H110").Value = 618.5
$ws.Range("I110").Value = 618.5
$ws.Range("K110").Value = 618.5
$ws.Range("M110").Value = 1426.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 27257.8
$ws.Range("I99").Value = 47349.637
$ws.Range("K99").Value = 47349.637
$ws.Range("M99").Value = -45851.637
$ws.Range("H105").Value = 1453.3478
$ws.Range("I105").Value = 1474
$ws.Range("K105").Value = 1474
$ws.Range("M105").Value = 273
$ws.Range("H107").Value = 5097.5454
$ws.Range("I107").Value = 4607.4
$ws.Range("J107").Value = 9999
$ws.Range("K107").Value = 4607.4
$ws.Range("L107").Value = 9999
$ws.Range("M107").Value = -2687.4
$ws.Range("N107").Value = -13839
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H37").Value = 3000
$ws.Range("J37").Value = 3000
$ws.Range("L37").Value = 3000
$ws.Range("N37").Value = -3214
$ws.Range("H58").Value = 590528.25
$ws.Range("I58").Value = 824677.7
$ws.Range("J58").Value = 5154.6665
$ws.Range("K58").Value = 824677.7
$ws.Range("L58").Value = 5154.6665
$ws.Range("M58").Value = -824474.7
$ws.Range("N58").Value = -5560.6665
$ws.Range("H132").Value = 27818506
$ws.Range("I132").Value = 51651.145
$ws.Range("J132").Value = 125002500
$ws.Range("K132").Value = 154953.435
$ws.Range("L132").Value = 375007500
$ws.Range("M132").Value = -152423.435
$ws.Range("N132").Value = -375012560
$ws.Range("H136").Value = 590528.25
$ws.Range("I136").Value = 824677.7
$ws.Range("J136").Value = 5154.6665
$ws.Range("K136").Value = 2474033.1
$ws.Range("L136").Value = 15463.9995
$ws.Range("M136").Value = -2471483.1
$ws.Range("N136").Value = -20563.9995
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3327
$ws.Range("I122").Value = 3124
$ws.Range("K122").Value = 9372
$ws.Range("M122").Value = -6922
$ws.Range("H126").Value = 2382484
$ws.Range("I126").Value = 8334444
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 25003332
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -25000862
$ws.Range("N126").Value = -10040
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2884.5454
$ws.Range("I7").Value = 2691.25
$ws.Range("K7").Value = 2691.25
$ws.Range("M7").Value = -2579.25
$ws.Range("H22").Value = 103162.7
$ws.Range("I22").Value = 169275
$ws.Range("K22").Value = 169275
$ws.Range("M22").Value = -168980
$ws.Range("H27").Value = 103162.7
$ws.Range("I27").Value = 169275
$ws.Range("K27").Value = 169275
$ws.Range("M27").Value = -169168
$ws.Range("H43").Value = 602400
$ws.Range("J43").Value = 1000000
$ws.Range("L43").Value = 1000000
$ws.Range("N43").Value = -1000386
$ws.Range("H55").Value = 1129
$ws.Range("I55").Value = 205.16667
$ws.Range("K55").Value = 205.16667
$ws.Range("M55").Value = -32.16667000000001
$ws.Range("H82").Value = 2658.9285
$ws.Range("I82").Value = 1539.8572
$ws.Range("J82").Value = 3778
$ws.Range("K82").Value = 1539.8572
$ws.Range("L82").Value = 3778
$ws.Range("M82").Value = -1178.8572
$ws.Range("N82").Value = -4500
$ws.Range("H85").Value = 2658.9285
$ws.Range("I85").Value = 1539.8572
$ws.Range("J85").Value = 3778
$ws.Range("K85").Value = 1539.8572
$ws.Range("L85").Value = 3778
$ws.Range("M85").Value = -291.8571999999999
$ws.Range("N85").Value = -6274
$ws.Range("H93").Value = 1452
$ws.Range("J93").Value = 2014.7273
$ws.Range("L93").Value = 2014.7273
$ws.Range("N93").Value = -4510.7273
$ws.Range("H122").Value = 3688.3416
$ws.Range("I122").Value = 3504.8462
$ws.Range("K122").Value = 10514.5386
$ws.Range("M122").Value = -8064.5386
$ws.Range("H126").Value = 2884.5454
$ws.Range("I126").Value = 2691.25
$ws.Range("K126").Value = 8073.75
$ws.Range("M126").Value = -5603.75
$ws.Range("H136").Value = 11999.667
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").Value = $null
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2674.2122
$ws.Range("I113").Value = 1369.75
$ws.Range("K113").Value = 4109.25
$ws.Range("M113").Value = -1939.25
$ws.Range("H122").Value = 2146.4075
$ws.Range("I122").Value = 1622.7
$ws.Range("J122").Value = 3642.7144
$ws.Range("K122").Value = 4868.1
$ws.Range("L122").Value = 10928.1432
$ws.Range("M122").Value = -2418.1
$ws.Range("N122").Value = -15828.1432
$ws.Range("H126").Value = 3616.7058
$ws.Range("I126").Value = 2707.25
$ws.Range("J126").Value = 5799.4
$ws.Range("K126").Value = 8121.75
$ws.Range("L126").Value = 17398.2
$ws.Range("M126").Value = -5651.75
$ws.Range("N126").Value = -22338.2
